$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 6

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G3"))

$ws.Range("G1").Value = "total_clp"
$ws.Range("G2").Value = 20083803
$ws.Range("G3").Value = 20083803

$ws.Columns.Item(7).ColumnWidth = 10.67

$ws.Range("J4").Select()
